# Updates cached market-data values across all item sheets (ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR) per the latest scheduled market-board refresh.
# Each entry below is a worksheet name -> list of {Cell, Value} pairs to
# write; a $null Value means the cell is cleared (no profit could be
# computed for that leve on this pass).

$wb = $excel.ActiveWorkbook

$changes = @{
    "ALC" = @(
        @{Cell="H17"; Value=911.40424},
        @{Cell="J17"; Value=900.8182},
        @{Cell="L17"; Value=2702.4546},
        @{Cell="N17"; Value=-3038.4546},
        @{Cell="H40"; Value=4121.4},
        @{Cell="J40"; Value=4409.074},
        @{Cell="L40"; Value=4409.074},
        @{Cell="N40"; Value=-4759.074},
        @{Cell="H41"; Value=869.6818},
        @{Cell="I41"; Value=865.6667},
        @{Cell="J41"; Value=872.46155},
        @{Cell="K41"; Value=865.6667},
        @{Cell="L41"; Value=872.46155},
        @{Cell="M41"; Value=-425.6667},
        @{Cell="N41"; Value=-1752.46155},
        @{Cell="H51"; Value=6010.067},
        @{Cell="I51"; Value=2657.1428},
        @{Cell="J51"; Value=7030.522},
        @{Cell="K51"; Value=2657.1428},
        @{Cell="L51"; Value=7030.522},
        @{Cell="M51"; Value=-2173.1428},
        @{Cell="N51"; Value=-7998.522},
        @{Cell="H70"; Value=5828.551},
        @{Cell="J70"; Value=6286.2046},
        @{Cell="L70"; Value=18858.6138},
        @{Cell="N70"; Value=-19398.6138},
        @{Cell="H73"; Value=5828.551},
        @{Cell="J73"; Value=6286.2046},
        @{Cell="L73"; Value=18858.6138},
        @{Cell="N73"; Value=-20730.6138},
        @{Cell="H96"; Value=643.55554},
        @{Cell="I96"; Value=502.42856},
        @{Cell="K96"; Value=1507.28568},
        @{Cell="M96"; Value=-134.28568},
        @{Cell="H100"; Value=3120.8},
        @{Cell="I100"; Value=3023.2222},
        @{Cell="K100"; Value=3023.2222},
        @{Cell="M100"; Value=-2482.2222},
        @{Cell="H133"; Value=65999},
        @{Cell="I133"; Value=65999},
        @{Cell="K133"; Value=65999},
        @{Cell="M133"; Value=-60939},
        @{Cell="H138"; Value=1982.92},
        @{Cell="I138"; Value=958.7727},
        @{Cell="J138"; Value=2787.6072},
        @{Cell="K138"; Value=2876.3181},
        @{Cell="L138"; Value=8362.821599999999},
        @{Cell="M138"; Value=2263.6819},
        @{Cell="N138"; Value=-18642.8216}
    )
    "ARM" = @(
        @{Cell="H16"; Value=3262.6667},
        @{Cell="I16"; Value=1444},
        @{Cell="J16"; Value=6900},
        @{Cell="K16"; Value=1444},
        @{Cell="L16"; Value=6900},
        @{Cell="M16"; Value=-1157},
        @{Cell="N16"; Value=-7474},
        @{Cell="H34"; Value=25},
        @{Cell="I34"; Value=25},
        @{Cell="K34"; Value=25},
        @{Cell="M34"; Value=246},
        @{Cell="I74"; Value=109727.09},
        @{Cell="J74"; Value=296643},
        @{Cell="K74"; Value=109727.09},
        @{Cell="L74"; Value=296643},
        @{Cell="M74"; Value=-108853.09},
        @{Cell="N74"; Value=-298391},
        @{Cell="I77"; Value=109727.09},
        @{Cell="J77"; Value=296643},
        @{Cell="K77"; Value=548635.45},
        @{Cell="L77"; Value=1483215},
        @{Cell="M77"; Value=-544267.45},
        @{Cell="N77"; Value=-1491951},
        @{Cell="H94"; Value=20082},
        @{Cell="J94"; Value=20082},
        @{Cell="L94"; Value=20082},
        @{Cell="N94"; Value=-21884},
        @{Cell="H122"; Value=746179.75},
        @{Cell="I122"; Value=2096.8572},
        @{Cell="K122"; Value=6290.571599999999},
        @{Cell="M122"; Value=-3840.571599999999}
    )
    "BSM" = @(
        @{Cell="H20"; Value=2610.4},
        @{Cell="I20"; Value=2224.6},
        @{Cell="J20"; Value=3382},
        @{Cell="K20"; Value=2224.6},
        @{Cell="L20"; Value=3382},
        @{Cell="M20"; Value=-1977.6},
        @{Cell="N20"; Value=-3876},
        @{Cell="H64"; Value=2159.4},
        @{Cell="I64"; Value=2099},
        @{Cell="J64"; Value=2174.5},
        @{Cell="K64"; Value=2099},
        @{Cell="L64"; Value=2174.5},
        @{Cell="M64"; Value=-1874},
        @{Cell="N64"; Value=-2624.5},
        @{Cell="H67"; Value=2159.4},
        @{Cell="I67"; Value=2099},
        @{Cell="J67"; Value=2174.5},
        @{Cell="K67"; Value=2099},
        @{Cell="L67"; Value=2174.5},
        @{Cell="M67"; Value=-1319},
        @{Cell="N67"; Value=-3734.5},
        @{Cell="H86"; Value=3576604.2},
        @{Cell="I86"; Value=5562362.5},
        @{Cell="K86"; Value=5562362.5},
        @{Cell="M86"; Value=-5561239.5},
        @{Cell="H89"; Value=3576604.2},
        @{Cell="I89"; Value=5562362.5},
        @{Cell="K89"; Value=27811812.5},
        @{Cell="M89"; Value=-27806196.5},
        @{Cell="H134"; Value=6291.6665},
        @{Cell="I134"; Value=1756.6},
        @{Cell="J134"; Value=8559.200000000001},
        @{Cell="K134"; Value=5269.799999999999},
        @{Cell="L134"; Value=25677.6},
        @{Cell="M134"; Value=-2734.799999999999},
        @{Cell="N134"; Value=-30747.6}
    )
    "CRP" = @(
        @{Cell="H50"; Value=1985.25},
        @{Cell="J50"; Value=1985.25},
        @{Cell="L50"; Value=1985.25},
        @{Cell="N50"; Value=-3235.25},
        @{Cell="H132"; Value=120353.75},
        @{Cell="I132"; Value=78675.16},
        @{Cell="K132"; Value=236025.48},
        @{Cell="M132"; Value=-233495.48}
    )
    "CUL" = @(
        @{Cell="H34"; Value=1020},
        @{Cell="I34"; Value=140.66667},
        @{Cell="J34"; Value=1899.3334},
        @{Cell="K34"; Value=422.00001},
        @{Cell="L34"; Value=5698.0002},
        @{Cell="M34"; Value=-338.00001},
        @{Cell="N34"; Value=-5866.0002},
        @{Cell="H39"; Value=2398.1428},
        @{Cell="I39"; Value=1500},
        @{Cell="J39"; Value=2547.8333},
        @{Cell="K39"; Value=4500},
        @{Cell="L39"; Value=7643.499899999999},
        @{Cell="M39"; Value=-4206},
        @{Cell="N39"; Value=-8231.499899999999},
        @{Cell="H55"; Value=55350.79},
        @{Cell="I55"; Value=708.75},
        @{Cell="J55"; Value=95090.45},
        @{Cell="K55"; Value=2126.25},
        @{Cell="L55"; Value=285271.35},
        @{Cell="M55"; Value=-1949.25},
        @{Cell="N55"; Value=-285625.35},
        @{Cell="H86"; Value=408.33334},
        @{Cell="I86"; Value=376.66666},
        @{Cell="J86"; Value=440},
        @{Cell="K86"; Value=1129.99998},
        @{Cell="L86"; Value=1320},
        @{Cell="M86"; Value=56.00001999999995},
        @{Cell="N86"; Value=-3692},
        @{Cell="H89"; Value=408.33334},
        @{Cell="I89"; Value=376.66666},
        @{Cell="J89"; Value=440},
        @{Cell="K89"; Value=3389.99994},
        @{Cell="L89"; Value=3960},
        @{Cell="M89"; Value=2538.00006},
        @{Cell="N89"; Value=-15816},
        @{Cell="H107"; Value=289.97058},
        @{Cell="I107"; Value=187.88},
        @{Cell="J107"; Value=573.55554},
        @{Cell="K107"; Value=563.64},
        @{Cell="L107"; Value=1720.66662},
        @{Cell="M107"; Value=1356.36},
        @{Cell="N107"; Value=-5560.66662},
        @{Cell="H113"; Value=2185.5454},
        @{Cell="J113"; Value=1474.5333},
        @{Cell="L113"; Value=4423.5999},
        @{Cell="N113"; Value=-8763.599900000001},
        @{Cell="H131"; Value=13442966},
        @{Cell="I131"; Value=11905577},
        @{Cell="K131"; Value=35716731},
        @{Cell="M131"; Value=-35711691},
        @{Cell="H139"; Value=71430830},
        @{Cell="I139"; Value=125001820},
        @{Cell="K139"; Value=375005460},
        @{Cell="M139"; Value=-375000320}
    )
    "GSM" = @(
        @{Cell="H49"; Value=50017500},
        @{Cell="I49"; Value=100000000},
        @{Cell="K49"; Value=100000000},
        @{Cell="M49"; Value=-99999816},
        @{Cell="H52"; Value=45000},
        @{Cell="J52"; Value=45000},
        @{Cell="L52"; Value=45000},
        @{Cell="N52"; Value=-45518},
        @{Cell="H70"; Value=8005545},
        @{Cell="I70"; Value=11116374},
        @{Cell="K70"; Value=11116374},
        @{Cell="M70"; Value=-11116104},
        @{Cell="H73"; Value=8005545},
        @{Cell="I73"; Value=11116374},
        @{Cell="K73"; Value=11116374},
        @{Cell="M73"; Value=-11115438},
        @{Cell="H122"; Value=237070.5},
        @{Cell="I122"; Value=298423.5},
        @{Cell="K122"; Value=895270.5},
        @{Cell="M122"; Value=-892820.5},
        @{Cell="H132"; Value=3901.4736},
        @{Cell="I132"; Value=3348.4614},
        @{Cell="K132"; Value=10045.3842},
        @{Cell="M132"; Value=-7515.3842}
    )
    "LTW" = @(
        @{Cell="H29"; Value=0},
        @{Cell="J29"; Value=0},
        @{Cell="L29"; Value=0},
        @{Cell="N29"; Value=$null},
        @{Cell="H41"; Value=45000},
        @{Cell="I41"; Value=0},
        @{Cell="K41"; Value=0},
        @{Cell="M41"; Value=$null},
        @{Cell="H55"; Value=1906.359},
        @{Cell="I55"; Value=1488.8214},
        @{Cell="J55"; Value=2969.182},
        @{Cell="K55"; Value=1488.8214},
        @{Cell="L55"; Value=2969.182},
        @{Cell="M55"; Value=-1315.8214},
        @{Cell="N55"; Value=-3315.182},
        @{Cell="H68"; Value=2521.6924},
        @{Cell="J68"; Value=2582.5},
        @{Cell="L68"; Value=2582.5},
        @{Cell="N68"; Value=-4080.5},
        @{Cell="H71"; Value=2521.6924},
        @{Cell="J71"; Value=2582.5},
        @{Cell="L71"; Value=12912.5},
        @{Cell="N71"; Value=-20400.5}
    )
    "WVR" = @(
        @{Cell="H113"; Value=1038.8928},
        @{Cell="I113"; Value=783.9474},
        @{Cell="K113"; Value=2351.8422},
        @{Cell="M113"; Value=-181.8422},
        @{Cell="H122"; Value=2524.5881},
        @{Cell="I122"; Value=1785.25},
        @{Cell="K122"; Value=5355.75},
        @{Cell="M122"; Value=-2905.75},
        @{Cell="H135"; Value=123994.5},
        @{Cell="J135"; Value=149990},
        @{Cell="L135"; Value=149990},
        @{Cell="N135"; Value=-160130}
    )
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $changes[$sheetName]) {
        $cell = $ws.Range($entry.Cell)
        if ($null -eq $entry.Value) {
            $cell.ClearContents()
        } else {
            $cell.Value = $entry.Value
        }
    }
}
